# Updates work week Feb 9th 2026
# Add a new "Loki" row (Grafana Loki migration guide link) to the links table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data for row 13
$newTitle = "Loki"
$newUrl   = "https://rewe.atlassian.net/wiki/spaces/O11Y/pages/1211416349/User+Guide+Logscale+Grafana+Loki+Migration"

# Write the URL first, then the title (matches the shared-string insertion
# order used when this row was originally authored)
$ws.Range("B13").Value = $newUrl
$ws.Range("A13").Value = $newTitle

# Add the actual hyperlink on B13 (this also applies the "Hyperlink" cell
# style, matching the other URL cells such as B12)
$ws.Hyperlinks.Add($ws.Range("B13"), $newUrl, [Type]::Missing, [Type]::Missing, $newUrl) | Out-Null

# Resize / refresh the Excel table (ListObject) to include the new row
$tbl = $ws.ListObjects.Item("Table3")
$tbl.Resize($ws.Range("A3:C13"))

# Update view: scroll and select B13 as the active cell
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B13").Select() | Out-Null

$wb.Save()
